$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the FAQ content (values replace the old Q/A pairs; three new rows
#    are appended for the additional FAQ entries).
# ---------------------------------------------------------------------------

$ws.Range("A2").Value = "What are the required files in uploading a Payment Screenshot?"
$ws.Range("B2").Value = "We accept JPG, JPEG, PNG, and HEIC file formats for image uploads."

$ws.Range("A3").Value = "What time is your check in and check out?"
$ws.Range("B3").Value = "Check-in time is at 4:00 PM and Check-out is at 12:00 NN."

$ws.Range("A4").Value = "Are pets allowed during our stay in your resort?"
$ws.Range("B4").Value = "Yes, pets are allowed but owners must be responsible for cleaning up after them."

$ws.Range("A5").Value = "Are there any restrictions on bringing outside food and drinks?"
$ws.Range("B5").Value = "There are no restrictions on bringing outside food and drinks."

$ws.Range("A6").Value = "Does resort offer wifi? Is it free or there are additional fee?"
$ws.Range("B6").Value = "Yes, we offer free WiFi for guests who have completed their booking, though coverage may be limited in some areas."

$ws.Range("A7").Value = "What are the payment methods? Do you accept full cash?"
$ws.Range("B7").Value = "We require 50% initial payment through online wallets like GCash. The remaining balance can be settled in cash upon arrival."

# ---------------------------------------------------------------------------
# 2. Column A grows wider to fit the longer questions.
# ---------------------------------------------------------------------------

$ws.Columns.Item(1).ColumnWidth = 51.63

# ---------------------------------------------------------------------------
# 3. Formatting. A2 becomes a wrapped 11pt dark-gray cell; B2 gets the same
#    font plus a white fill; the rest of the answer rows (3-7) switch to the
#    11pt "theme" colored font; trailing blank rows (8-19) are pre-formatted
#    with that same font so the sheet has a consistent look further down.
# ---------------------------------------------------------------------------

$ws.Range("A2").Font.Size = 11
$ws.Range("A2").Font.Color = 4210752
$ws.Range("A2").Font.Name = "Arial"
$ws.Range("A2").WrapText = $true

$ws.Range("B2").Font.Size = 11
$ws.Range("B2").Font.Color = 4210752
$ws.Range("B2").Font.Name = "Arial"
$ws.Range("B2").Interior.Color = 16777215

$ws.Range("A3").Font.Size = 11
$ws.Range("A3").Font.Name = "Arial"

$ws.Range("A3").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("B7").PasteSpecial(-4122)

# Pre-format 12 trailing blank rows beneath the table (rows 8-19).
$ws.Range("A8").Font.Size = 11
$ws.Range("A8").Font.Name = "Arial"

$ws.Range("A8").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("B19").PasteSpecial(-4122)
